$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: recovery-factor parameter change + new R^2/ABS column ---
$ws.Range("K2").Value = 1
$ws.Range("G2").Formula = '=ABS(D2-E2)'

# --- Input (B column) corrections on existing rows ---
$ws.Range("B6").Value = 1000
$ws.Range("B8").Value = 1000
$ws.Range("B10").Value = 1000
$ws.Range("B11").Value = 1000

# --- D column: record (hard-code) the computed outputs on rows whose
#     inputs changed or which are no longer driven live by the formula ---
$ws.Range("D4").Value = 100
$ws.Range("D6").Value = 100
$ws.Range("D8").Value = 100
$ws.Range("D9").Value = 100
$ws.Range("D10").Value = 100
$ws.Range("D11").Value = 100

# --- New G column (ABS deviation) for existing rows 3-11 ---
$ws.Range("G3:G11").Formula = '=ABS(D3-E3)'

# --- New H column (recovery factor helper) for existing rows ---
$ws.Range("H3").Formula = '=(G2^2+G3^2)/(G2+G3)'
$ws.Range("H4").Formula = '=(G4+G3)/2'
$ws.Range("H5:H11").Formula = '=(G5+G4)/2'

# --- Append new data rows 12-19 (A=11..18, constant inputs) ---
# NOTE: row 12 previously held the old totals formulas (D12/E12/F12); writing
# the new literal/formula content directly over those cells replaces them
# cleanly (doing a separate ClearContents afterwards would corrupt the
# F12:F19 shared-formula group created below).
$newRows = @(
    @{ Row = 12; A = 11 },
    @{ Row = 13; A = 12 },
    @{ Row = 14; A = 13 },
    @{ Row = 15; A = 14 },
    @{ Row = 16; A = 15 },
    @{ Row = 17; A = 16 },
    @{ Row = 18; A = 17 },
    @{ Row = 19; A = 18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = 1000
    $ws.Cells.Item($row, 3).Value = 0.1
    $ws.Cells.Item($row, 4).Value = 100
}

# E column formulas for the new rows (ROUNDUP growing-range formula)
$ws.Range("E12").Formula = '=ROUNDUP((SUM(D$2:D11) - SUM($E$2:E11))*$K$2 + D12, 0)'
$ws.Range("E13").Formula = '=ROUNDUP((SUM(D$2:D12) - SUM($E$2:E12))*$K$2 + D13, 0)'
$ws.Range("E14").Formula = '=ROUNDUP((SUM(D$2:D13) - SUM($E$2:E13))*$K$2 + D14, 0)'
$ws.Range("E15").Formula = '=ROUNDUP((SUM(D$2:D14) - SUM($E$2:E14))*$K$2 + D15, 0)'
$ws.Range("E16").Formula = '=ROUNDUP((SUM(D$2:D15) - SUM($E$2:E15))*$K$2 + D16, 0)'
$ws.Range("E17").Formula = '=ROUNDUP((SUM(D$2:D16) - SUM($E$2:E16))*$K$2 + D17, 0)'
$ws.Range("E18").Formula = '=ROUNDUP((SUM(D$2:D17) - SUM($E$2:E17))*$K$2 + D18, 0)'
$ws.Range("E19").Formula = '=ROUNDUP((SUM(D$2:D18) - SUM($E$2:E18))*$K$2 + D19, 0)'

# F, G, H columns for the new rows (same relative pattern as existing rows).
# These writes directly replace the old SUM formulas that lived in D12/E12/F12.
$ws.Range("F12:F19").Formula = '=(D12-E12)^2'
$ws.Range("G12:G19").Formula = '=ABS(D12-E12)'
$ws.Range("H12:H19").Formula = '=(G12+G11)/2'

# --- New totals row at 20 (replaces the old row-12 totals) ---
$ws.Range("D20").Formula = '=SUM(D2:D13)'
$ws.Range("E20").Formula = '=SUM(E2:E13)'
$ws.Range("F20").Formula = '=SUM(F2:F13)'
$ws.Range("H20").Formula = '=SUM(H3:H13)'

# --- Selection bookkeeping (matches authored workbook state) ---
$ws.Range("D2").Select()

# --- Chart source ranges now cover the extended data (rows 2-19) ---
$chart = $ws.ChartObjects(1).Chart
$series1 = $chart.SeriesCollection(1)
$series1.Formula = "=SERIES(Sheet1!`$D`$1,Sheet1!`$A`$2:`$A`$19,Sheet1!`$D`$2:`$D`$19,1)"
$series2 = $chart.SeriesCollection(2)
$series2.Formula = "=SERIES(Sheet1!`$E`$1,Sheet1!`$A`$2:`$A`$19,Sheet1!`$E`$2:`$E`$19,2)"

# --- Reposition / resize the chart on the sheet ---
# (columns default to 58.4375pt, rows to 16pt on this sheet, so these
# Left/Top/Width/Height values reproduce the authored anchor exactly:
# from col 8 / row 5  ->  to col 18 / row 40)
$co = $ws.ChartObjects(1)
$co.Left = 489.5
$co.Top = 85
$co.Width = 609.375
$co.Height = 565
